# Apply "Meeting Recording & Burndown Update" changes to the burndown sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated actual-hours-spent figures (column D = "Day 2") for several tasks,
# plus revised estimated-hours (column B) for a few tasks that were re-scoped.

# Basement Room B: Assets & Props (row 9)
$ws.Range("D9").Value = 2

# Basement Room C: Assets & Props (row 10)
$ws.Range("D10").Value = 1

# Basement Room E: Assets & Props (row 11)
$ws.Range("D11").Value = 1.5

# Room Allocation System (row 13)
$ws.Range("D13").Value = 4

# Basement Room A: Event (row 15)
$ws.Range("D15").Value = 4

# Basement Room B: Event (row 16)
$ws.Range("B16").Value = 2
$ws.Range("D16").Value = 1

# Basement Room D: Event (row 17)
$ws.Range("B17").Value = 2

# User Stories (Final Event) (row 21)
$ws.Range("B21").Value = 2
$ws.Range("D21").Value = 0.2

# Reflect the cell the user last had selected
$ws.Range("D17").Select()

$excel.CalculateFullRebuild()
$excel.Calculate()

$wb.Save()
